$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# A2: 84711 -> 84712
$ws.Range("A2").Value = 84712

# P2: update locality name
$ws.Range("P2").Value = "Skäftekärr, Tujabeståndet, Öl"

# S2: accuracy 25 -> 10
$ws.Range("S2").Value = 10

# X2: new external id value (previously empty, inserted between W2 and Y2)
$ws.Range("X2").Value = "Hö-Bor-1880"

# AI2: biotope description removed entirely
$ws.Range("AI2").ClearContents()

# AW2: reporter changed
$ws.Range("AW2").Value = "Öland- Floraväktarna"

# AX2: observers list trimmed
$ws.Range("AX2").Value = "Thomas Gunnarsson, Ulla-Britt Andersson"

# AY2: project name added
$ws.Range("AY2").Value = "Floraväkteri Sverige"
